$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Final Extraction Sheet")
$lo = $ws.ListObjects.Item(1)
$ws.Columns.Item(2).Insert()
$lo.Resize($ws.Range("A1:Q15"))
for ($i=1; $i -le $lo.ListColumns.Count; $i++) {
    Write-Output "$i : $($lo.ListColumns.Item($i).Name)"
}
